$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two more matches are appended for Kieron Pollard (c):
#   row 4 - Abu Dhabi, 25 Oct 2020 vs Rajasthan Royals  (same stats as row 3)
#   row 5 - Abu Dhabi, 28 Oct 2020 vs Royal Challengers Bangalore (same stats as row 2)
# Copy the existing rows so the numeric-looking columns (G:K) keep being
# stored as text, exactly like the rest of the sheet, instead of being
# auto-converted to real numbers.
$ws.Range("A3:K3").Copy()
$ws.Range("A4:K4").PasteSpecial()

$ws.Range("A2:K2").Copy()
$ws.Range("A5:K5").PasteSpecial()

$excel.CutCopyMode = $false
